$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 73347.2945
$ws.Range("H2").Value = 2723.589000000007

$ws.Range("A3").Value = 70244.77099999999
$ws.Range("H3").Value = 19786.54199999999

$ws.Range("A4").Value = 67129.0255
$ws.Range("H4").Value = 16621.05100000001

$ws.Range("A5").Value = 64595.502
$ws.Range("H5").Value = 14120.004

$ws.Range("A6").Value = 65374.48
$ws.Range("H6").Value = 14904.96000000001

$ws.Range("A7").Value = 67982.29700000001
$ws.Range("H7").Value = 17406.594

$ws.Range("A8").Value = 67367.1265
$ws.Range("H8").Value = 38725.253

$ws.Range("A9").Value = 79207.3645
$ws.Range("H9").Value = 27043.72899999999

$ws.Range("A10").Value = 95174.5395
$ws.Range("H10").Value = 21527.079

$ws.Range("A11").Value = 81952.2095
$ws.Range("H11").Value = 4516.418999999994

$ws.Range("A12").Value = 81626.978
$ws.Range("H12").Value = 1155.956000000006

$ws.Range("A13").Value = 80925.1525
$ws.Range("H13").Value = -459.695000000007

$ws.Range("A14").Value = 80952.772
$ws.Range("H14").Value = -109.4560000000056

$ws.Range("A15").Value = 107586.796
$ws.Range("H15").Value = 26401.592

$ws.Range("A16").Value = 106447.456
$ws.Range("H16").Value = 24921.91200000001

$ws.Range("A17").Value = 65630.38250000001
$ws.Range("H17").Value = -15832.23499999999

$ws.Range("A18").Value = 40052.2725
$ws.Range("H18").Value = -41709.455

$ws.Range("H19").Value = 7051.228499999997
$ws.Range("H20").Value = 8134.193499999994
$ws.Range("H21").Value = 8380.1875
$ws.Range("H22").Value = 8944.044999999998
$ws.Range("H23").Value = 10285.514
$ws.Range("H24").Value = 10044.24250000001
$ws.Range("H25").Value = -18483.906
